$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column-D target cells to Text format so numeric-looking values
# (e.g. 297.92, 0.0693) are stored as literal text, matching the source data.
$dCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D9', 'D10', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D21', 'D22', 'D23', 'D24', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D37', 'D39', 'D43', 'D44', 'D46', 'D47', 'D50', 'D51')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '42.048.71'
$ws.Range("E2").Value = '  -2.02%  '

# Row 3
$ws.Range("D3").Value = '2.258.94'
$ws.Range("E3").Value = '  -3.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '297.92'
$ws.Range("E5").Value = '  -2.93%  '

# Row 6
$ws.Range("D6").Value = '93.62'
$ws.Range("E6").Value = '  -7.01%  '

# Row 7
$ws.Range("D7").Value = '0.498'
$ws.Range("E7").Value = '  -2.33%  '

# Row 8
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -3.86%  '

# Row 10
$ws.Range("D10").Value = '32.93'
$ws.Range("E10").Value = '  -5.84%  '

# Row 11
$ws.Range("E11").Value = '  -1.77%  '

# Row 12
$ws.Range("D12").Value = '47.71'
$ws.Range("E12").Value = '  -8.61%  '

# Row 13
$ws.Range("D13").Value = '0.113'
$ws.Range("E13").Value = '  +0.51%  '

# Row 14
$ws.Range("D14").Value = '6.66'
$ws.Range("E14").Value = '  -2.56%  '

# Row 15
$ws.Range("D15").Value = '2.611.13'
$ws.Range("E15").Value = '  -3.53%  '

# Row 16
$ws.Range("D16").Value = '15.28'
$ws.Range("E16").Value = '  -3.71%  '

# Row 17
$ws.Range("D17").Value = '2.260.14'
$ws.Range("E17").Value = '  -1.91%  '

# Row 18
$ws.Range("D18").Value = '0.775'
$ws.Range("E18").Value = '  -3.52%  '

# Row 19
$ws.Range("D19").Value = '42.064.72'
$ws.Range("E19").Value = '  -1.78%  '

# Row 20
$ws.Range("E20").Value = '  -2.48%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.01'
$ws.Range("E21").Value = '  -3.63%  '

# Row 22
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").Value = '11.37'
$ws.Range("E22").Value = '  -3.01%  '

# Row 23
$ws.Range("D23").Value = '66.55'
$ws.Range("E23").Value = '  -2.01%  '

# Row 24
$ws.Range("D24").Value = '232.97'
$ws.Range("E24").Value = '  -1.72%  '

# Row 25
$ws.Range("E25").Value = '  -4.14%  '

# Row 26
$ws.Range("E26").Value = '  +0.09%  '

# Row 27
$ws.Range("D27").Value = '2.45'
$ws.Range("E27").Value = '  -4.41%  '

# Row 28
$ws.Range("D28").Value = '23.70'
$ws.Range("E28").Value = '  -7.33%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.16'
$ws.Range("E29").Value = '  -7.18%  '

# Row 30
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '166.42'
$ws.Range("E30").Value = '  +4.03%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '33.53'
$ws.Range("E31").Value = '  -4.21%  '

# Row 32
$ws.Range("D32").Value = '9.03'
$ws.Range("E32").Value = '  -3.61%  '

# Row 33
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("D34").Value = '4.92'
$ws.Range("E34").Value = '  -4.11%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '2.35'
$ws.Range("E35").Value = '  -4.45%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.0693'
$ws.Range("E36").Value = '  -4.70%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.38'
$ws.Range("E37").Value = '  -6.10%  '

# Row 38
$ws.Range("E38").Value = '  -5.62%  '

# Row 39
$ws.Range("D39").Value = '15.92'
$ws.Range("E39").Value = '  -8.20%  '

# Row 40
$ws.Range("E40").Value = '  -4.78%  '

# Row 41
$ws.Range("E41").Value = '  -3.46%  '

# Row 42
$ws.Range("E42").Value = '  -8.39%  '

# Row 43
$ws.Range("D43").Value = '2.40'
$ws.Range("E43").Value = '  +1.86%  '

# Row 44
$ws.Range("D44").Value = '1.940.20'
$ws.Range("E44").Value = '  -4.33%  '

# Row 45
$ws.Range("E45").Value = '  -2.51%  '

# Row 46
$ws.Range("D46").Value = '17.40'
$ws.Range("E46").Value = '  -6.96%  '

# Row 47
$ws.Range("D47").Value = '9.56'
$ws.Range("E47").Value = '  -7.60%  '

# Row 48
$ws.Range("E48").Value = '  -5.54%  '

# Row 49
$ws.Range("E49").Value = '  -3.05%  '

# Row 50
$ws.Range("D50").Value = '2.484.31'

# Row 51
$ws.Range("D51").Value = '52.21'
$ws.Range("E51").Value = '  -7.25%  '
